# Update countries & provincias Spain
# Refresh the COVID dashboard snapshot: update the "last updated" timestamp,
# update case figures for a handful of countries, and re-rank two countries
# (Polinesia Francesa and Burkina Faso) that moved up past their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 03:15"

# --- Estados Unidos (row 4): updated totals ---
$ws.Range("B4").Value = 7894338
$ws.Range("C4").Value = 60418
$ws.Range("D4").Value = 5064200
$ws.Range("E4").Value = 2611496
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 904
$ws.Range("H4").Value = 218642

# --- Paraguay (row 68): updated totals ---
$ws.Range("B68").Value = 48275
$ws.Range("C68").Value = 959
$ws.Range("D68").Value = 30643
$ws.Range("E68").Value = 16587
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 33
$ws.Range("H68").Value = 1045

# --- Gabon (row 113): updated totals ---
$ws.Range("B113").Value = 8835
$ws.Range("C113").Value = 20
$ws.Range("D113").Value = 8189
$ws.Range("E113").Value = 592

# --- Polinesia Francesa jumps ahead of Principado de Andorra and Letonia ---
# Row 150 becomes Polinesia Francesa with its new, higher totals.
$ws.Range("A150").Value = "Polinesia Francesa"
$ws.Range("B150").Value = 2754
$ws.Range("C150").Value = 334
$ws.Range("D150").Value = 2019
$ws.Range("E150").Value = 725
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 10

# Row 151 becomes Principado de Andorra (previous row 150 data, unchanged).
$ws.Range("A151").Value = "Principado de Andorra"
$ws.Range("B151").Value = 2696
$ws.Range("C151").Value = 128
$ws.Range("D151").Value = 1814
$ws.Range("E151").Value = 827
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 55

# Row 152 becomes Letonia (previous row 151 data, unchanged).
$ws.Range("A152").Value = "Letonia"
$ws.Range("B152").Value = 2507
$ws.Range("C152").Value = 137
$ws.Range("D152").Value = 1322
$ws.Range("E152").Value = 1145
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 40

# --- Burkina Faso jumps ahead of Uruguay ---
# Row 157 becomes Burkina Faso with its new, higher totals.
$ws.Range("A157").Value = "Burkina Faso"
$ws.Range("B157").Value = 2254
$ws.Range("C157").Value = 13
$ws.Range("D157").Value = 1516
$ws.Range("E157").Value = 678
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 60

# Row 158 becomes Uruguay (previous row 157 data, unchanged).
$ws.Range("A158").Value = "Uruguay"
$ws.Range("B158").Value = 2251
$ws.Range("C158").Value = 25
$ws.Range("D158").Value = 1917
$ws.Range("E158").Value = 285
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 49

# --- Niger (row 167): updated totals ---
$ws.Range("D167").Value = 1123
$ws.Range("E167").Value = 9

# --- Islas Turcas y Caicos (row 173): updated totals ---
$ws.Range("D173").Value = 672
$ws.Range("E173").Value = 17

# --- Islas Caimanes (row 190): updated totals ---
$ws.Range("B190").Value = 220
$ws.Range("C190").Value = 6
$ws.Range("E190").Value = 8
